$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateRecipient")
$ws.Activate()

# Best-effort: the saved window width shrank (Excel window was resized before
# the save that produced this commit). Not all hosts expose this, so ignore
# failures.
try {
    $excel.ActiveWindow.Width = 17415
} catch {
}

# New row 4: TrialData / Recipient / 91
$ws.Range("A4").Value = "TrialData"
$ws.Range("B4").Value = "'Recipient"
$ws.Range("C4").Value = "'91"

# New row 5: TrialData / Recipient / 27
$ws.Range("A5").Value = "TrialData"
$ws.Range("B5").Value = "'Recipient"
$ws.Range("C5").Value = "'27"

$ws.Range("C6").Select()
